$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.345.18"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "3.514.35"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'610.10"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").Value = "'150.66"
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("D7").Value = "3.512.11"
$ws.Range("E7").Value = "  -1.51%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").Value = "'7.03"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").Value = "'0.424"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("D14").Value = "4.106.77"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").Value = "'31.82"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "3.510.37"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("D17").Value = "67.332.23"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "'6.39"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("D21").Value = "'442.97"
$ws.Range("E21").Value = "  -2.78%  "
$ws.Range("D22").Value = "'9.29"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("D24").Value = "'77.32"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("E25").Value = "  +8.89%  "
$ws.Range("D26").Value = "3.652.92"
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "'10.24"
$ws.Range("E28").Value = "  -3.78%  "
$ws.Range("D29").Value = "'8.33"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.54"
$ws.Range("E32").Value = "  -5.96%  "
$ws.Range("D33").Value = "'0.164"
$ws.Range("E33").Value = "  +2.46%  "
$ws.Range("D34").Value = "'25.81"
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("E35").Value = "  -1.60%  "
$ws.Range("D36").Value = "3.503.57"
$ws.Range("E36").Value = "  -1.75%  "
$ws.Range("E37").Value = "  -4.08%  "
$ws.Range("D38").Value = "'8.00"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'177.35"
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").Value = "'2.17"
$ws.Range("E42").Value = "  +3.65%  "
$ws.Range("D43").Value = "'0.0875"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("D45").Value = "'0.880"
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("D46").Value = "'45.02"
$ws.Range("E46").Value = "  -2.51%  "
$ws.Range("D47").Value = "'27.57"
$ws.Range("E47").Value = "  -5.78%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.61"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "'1.26"
$ws.Range("E49").Value = "  +3.84%  "
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("E51").Value = "  -1.90%  "
